$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data: a start time logged in column A, plus a (mistakenly
# textual) "Time Spent " label in C13, and the actual time-spent value in C14.
$ws.Range("A13").Value = "01/12/2021 3:00 P.M."
$ws.Range("C13").Value = "Time Spent "
$ws.Range("C14").Value = 0.083333333333333329

# Move the active selection to C13 to match the saved selection state.
$ws.Range("C13").Select()
